# Slide 3 ("Използвани библиотеки") holds a SmartArt (vertical bullet list)
# graphic frame listing the Python libraries used in the project. One of the
# bullets is an unfinished / placeholder entry ("- аst -") that the author
# removed from the diagram.
#
# This script locates that SmartArt graphic, finds the bullet node whose text
# is the leftover placeholder, and deletes it - mirroring what a user does in
# the PowerPoint UI by selecting that line of the SmartArt and pressing
# Delete/Backspace (which removes the whole bullet/node).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# Find the shape that hosts the SmartArt graphic on this slide.
$smartArtShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasSmartArt) {
        $smartArtShape = $shp
    }
}

$sa = $smartArtShape.SmartArt
$nodes = $sa.AllNodes

# Locate the placeholder/unfinished bullet ("- аst -") and remove it. Walk
# backwards since deleting can shift indices.
for ($i = $nodes.Count; $i -ge 1; $i--) {
    $n = $nodes.Item($i)
    $txt = $n.TextFrame2.TextRange.Text
    if ($txt.Trim() -eq "- аst -") {
        $n.Delete()
    }
}
